$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.984.77"
$ws.Cells.Item(2, 5).Value = "  -2.47%  "

$ws.Cells.Item(3, 4).Value = "3.517.28"
$ws.Cells.Item(3, 5).Value = "  -3.08%  "

$ws.Cells.Item(5, 4).Value = "'586.14"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.35%  "

$ws.Cells.Item(6, 4).Value = "'169.94"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -3.51%  "

$ws.Cells.Item(7, 5).Value = "  -1.91%  "

$ws.Cells.Item(8, 4).Value = "3.511.24"
$ws.Cells.Item(8, 5).Value = "  -3.04%  "

$ws.Cells.Item(9, 5).Value = "  +0.02%  "

$ws.Cells.Item(10, 5).Value = "  -4.26%  "

$ws.Cells.Item(11, 4).Value = "'6.78"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -1.23%  "

$ws.Cells.Item(12, 5).Value = "  -5.58%  "

$ws.Cells.Item(13, 4).Value = "'47.36"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.45%  "

$ws.Cells.Item(14, 4).Value = "'0.0000274"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.35%  "

$ws.Cells.Item(15, 4).Value = "4.084.95"
$ws.Cells.Item(15, 5).Value = "  -3.06%  "

$ws.Cells.Item(16, 4).Value = "'8.40"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -6.57%  "

$ws.Cells.Item(17, 4).Value = "'611.69"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  -9.27%  "

$ws.Cells.Item(18, 4).Value = "69.068.38"
$ws.Cells.Item(18, 5).Value = "  -2.41%  "

$ws.Cells.Item(19, 4).Value = "3.517.10"
$ws.Cells.Item(19, 5).Value = "  -3.24%  "

$ws.Cells.Item(20, 5).Value = "  -1.73%  "

$ws.Cells.Item(21, 4).Value = "'17.37"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -2.49%  "

$ws.Cells.Item(22, 4).Value = "'11.08"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.65%  "

$ws.Cells.Item(23, 4).Value = "'0.884"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -6.48%  "

$ws.Cells.Item(24, 4).Value = "'15.72"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -8.74%  "

$ws.Cells.Item(25, 4).Value = "'96.51"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -3.60%  "

$ws.Cells.Item(26, 4).Value = "'3.83"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.61%  "

$ws.Cells.Item(27, 4).Value = "'0.999"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.01%  "

$ws.Cells.Item(28, 5).Value = "  -6.68%  "

$ws.Cells.Item(29, 4).Value = "'9.19"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -6.45%  "

$ws.Cells.Item(30, 4).Value = "'32.59"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -5.83%  "

$ws.Cells.Item(31, 5).Value = "  -7.17%  "

$ws.Cells.Item(32, 5).Value = "  -5.24%  "

$ws.Cells.Item(33, 4).Value = "'1.31"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -4.97%  "

$ws.Cells.Item(34, 4).Value = "'6.89"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -8.99%  "

$ws.Cells.Item(35, 4).Value = "'614.14"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +6.50%  "

$ws.Cells.Item(36, 4).Value = "'10.73"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -3.33%  "

$ws.Cells.Item(37, 4).Value = "'3.47"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -12.70%  "

$ws.Cells.Item(38, 5).Value = "  -5.47%  "

$ws.Cells.Item(39, 4).Value = "'56.98"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -2.62%  "

$ws.Cells.Item(40, 5).Value = "  +0.03%  "

$ws.Cells.Item(41, 5).Value = "  -1.94%  "

$ws.Cells.Item(42, 2).Value = "Kaspa"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(42, 4).Value = "'0.135"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -3.70%  "

$ws.Cells.Item(43, 2).Value = "Maker"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(43, 4).Value = "3.386.88"
$ws.Cells.Item(43, 5).Value = "  -4.44%  "

$ws.Cells.Item(44, 5).Value = "  -5.74%  "

$ws.Cells.Item(45, 4).Value = "'32.76"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -4.88%  "

$ws.Cells.Item(46, 4).Value = "0.0₃0696"

$ws.Cells.Item(47, 5).Value = "  -5.85%  "

$ws.Cells.Item(48, 5).Value = "  -8.13%  "

$ws.Cells.Item(49, 5).Value = "  -4.36%  "

$ws.Cells.Item(50, 4).Value = "'134.05"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.54%  "

$ws.Cells.Item(51, 4).Value = "'5.52"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +10.08%  "
